$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumbering of the two DAAS (and similar Ports/COS/DAAS) invalidated the
# stale DAAS name (E), physical port (F), PUERTOLIBRE tag (G) and the
# parsed-interface list (I) for rows 32-47, so those values are cleared.
# Column H (the numeric pair list, e.g. ['30','0']) stays as-is.
$ws.Range("E32:G47").ClearContents()
$ws.Range("I32:I47").ClearContents()
